# user_account.xlsx: add an "accesstype" column (G) for every existing user,
# and append 3 new user rows (admin, yoda, andrey) with correct access types,
# matching "correct login depending on user's accesstype".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing rows (1-4): tag each as "customer" in the new column G ---
$ws.Range("G1").Value = "customer"
$ws.Range("G2").Value = "customer"
$ws.Range("G3").Value = "customer"
$ws.Range("G4").Value = "customer"

# --- row 5: admin user ---
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "admin"
$ws.Range("C5").Value = "admin"
$ws.Range("D5").Value = "admin@mail.ru"
$ws.Range("E5").Value = "admin"
$ws.Range("F5").Value = "admin"
$ws.Range("G5").Value = "admin"

# --- row 6: vendor user "yoda" ---
$ws.Range("A6").Value = 6
$ws.Range("B6").Value = "yoda"
$ws.Range("C6").Value = "123yoda"
$ws.Range("D6").Value = "yoda@gmail.com"
$ws.Range("E6").Value = "yoda"
$ws.Range("F6").Value = "yyy"
$ws.Range("G6").Value = "vendor"

# --- row 7: customer user "andrey" ---
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "andrey"
$ws.Range("C7").Value = "andrey11223344"
$ws.Range("D7").Value = "andrin@mail.ru"
$ws.Range("E7").Value = "Andrey"
$ws.Range("F7").Value = "Skvortsov"
$ws.Range("G7").Value = "customer"

# --- hyperlink the new email cells, same as D1, and restore the
#     workbook's existing hyperlink cell style (Hyperlinks.Add swaps in
#     its own style xf, so re-apply the named style afterwards) ---
[void]$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:admin@mail.ru")
$ws.Range("D5").Style = "Гиперссылка"

[void]$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:yoda@gmail.com")
$ws.Range("D6").Style = "Гиперссылка"

# --- match final selection shown in the edited workbook ---
[void]$ws.Range("G6").Select()
